$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.5979736666666666
$ws.Range("H2").Value = 1.793921
$ws.Range("I2").Value = 0.03342655292740804
$ws.Range("J2").Value = 0.03342655292740804
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4374626666666666
$ws.Range("N2").Value = 1.312388
$ws.Range("O2").Value = 0.13789916414942
$ws.Range("P2").Value = 0.13789916414942
$ws.Range("Q2").Value = 0.2615911548164444
$ws.Range("R2").Value = 2.354320393347999
$ws.Range("S2").Value = 0.004609493709085918
$ws.Range("T2").Value = 0.004609493709085918

# Row 3
$ws.Range("G3").Value = 0.5979736666666666
$ws.Range("H3").Value = 1.793921
$ws.Range("I3").Value = 0.03342655292740804
$ws.Range("J3").Value = 0.03342655292740804
$ws.Range("O3").Value = 0.584321423572861
$ws.Range("P3").Value = 0.584321423572861
$ws.Range("Q3").Value = 1.108442657497111
$ws.Range("R3").Value = 9.975983917474
$ws.Range("S3").Value = 0.01953185099167665
$ws.Range("T3").Value = 0.01953185099167665

# Row 4
$ws.Range("G4").Value = 0.5979736666666666
$ws.Range("H4").Value = 1.793921
$ws.Range("I4").Value = 0.03342655292740804
$ws.Range("J4").Value = 0.03342655292740804
$ws.Range("O4").Value = 0.2777794122777191
$ws.Range("P4").Value = 0.2777794122777191
$ws.Range("Q4").Value = 0.5269403748033332
$ws.Range("R4").Value = 4.74246337323
$ws.Range("S4").Value = 0.009285208226645476
$ws.Range("T4").Value = 0.009285208226645478

# Row 5
$ws.Range("I5").Value = 0.8874158839838097
$ws.Range("J5").Value = 0.8874158839838097
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4374626666666666
$ws.Range("N5").Value = 1.312388
$ws.Range("O5").Value = 0.13789916414942
$ws.Range("P5").Value = 0.13789916414942
$ws.Range("Q5").Value = 6.944782682136443
$ws.Range("R5").Value = 62.503044139228
$ws.Range("S5").Value = 0.122373908654286
$ws.Range("T5").Value = 0.122373908654286

# Row 6
$ws.Range("I6").Value = 0.8874158839838097
$ws.Range("J6").Value = 0.8874158839838097
$ws.Range("O6").Value = 0.584321423572861
$ws.Range("P6").Value = 0.584321423572861
$ws.Range("S6").Value = 0.5185361126305885
$ws.Range("T6").Value = 0.5185361126305885

# Row 7
$ws.Range("I7").Value = 0.8874158839838097
$ws.Range("J7").Value = 0.8874158839838097
$ws.Range("O7").Value = 0.2777794122777191
$ws.Range("P7").Value = 0.2777794122777191
$ws.Range("S7").Value = 0.2465058626989352
$ws.Range("T7").Value = 0.2465058626989352

# Row 8
$ws.Range("I8").Value = 0.07915756308878232
$ws.Range("J8").Value = 0.07915756308878232
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4374626666666666
$ws.Range("N8").Value = 1.312388
$ws.Range("O8").Value = 0.13789916414942
$ws.Range("P8").Value = 0.13789916414942
$ws.Range("Q8").Value = 0.6194751336106666
$ws.Range("R8").Value = 5.575276202496
$ws.Range("S8").Value = 0.01091576178604806
$ws.Range("T8").Value = 0.01091576178604806

# Row 9
$ws.Range("I9").Value = 0.07915756308878232
$ws.Range("J9").Value = 0.07915756308878232
$ws.Range("O9").Value = 0.584321423572861
$ws.Range("P9").Value = 0.584321423572861
$ws.Range("S9").Value = 0.04625345995059583
$ws.Range("T9").Value = 0.04625345995059583

# Row 10
$ws.Range("I10").Value = 0.07915756308878232
$ws.Range("J10").Value = 0.07915756308878232
$ws.Range("O10").Value = 0.2777794122777191
$ws.Range("P10").Value = 0.2777794122777191
$ws.Range("Q10").Value = 0.5269403748033332
$ws.Range("R10").Value = 4.74246337323
$ws.Range("S10").Value = 0.009285208226645476
$ws.Range("T10").Value = 0.009285208226645478
